$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.296.95'
$ws.Range('E2').Value = '  +5.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.918.79'
$ws.Range('E3').Value = '  +6.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.74'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5183'
$ws.Range('E7').Value = '  +4.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.13'
$ws.Range('E8').Value = '  +6.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3001'
$ws.Range('E9').Value = '  +7.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06898'
$ws.Range('E10').Value = '  +8.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.918.84'
$ws.Range('E11').Value = '  +6.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '17.56'
$ws.Range('E12').Value = '  +4.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07345'
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6919'
$ws.Range('E14').Value = '  +7.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.20'
$ws.Range('E15').Value = '  +7.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.960'
$ws.Range('E16').Value = '  +5.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.298.34'
$ws.Range('E17').Value = '  +5.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008285'
$ws.Range('E18').Value = '  +12.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9980'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.14'
$ws.Range('E20').Value = '  +7.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.163.49'
$ws.Range('E21').Value = '  +6.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9989'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.886'
$ws.Range('E23').Value = '  +5.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.792'
$ws.Range('E24').Value = '  +9.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.232'
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.41'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('B27').Value = 'BitcoinCash'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '139.75'
$ws.Range('E27').Value = '  +25.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.37'
$ws.Range('E28').Value = '  +8.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.019'
$ws.Range('E29').Value = '  +7.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.382'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.307'
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08889'
$ws.Range('E32').Value = '  +6.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.038'
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05145'
$ws.Range('E34').Value = '  +3.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.166'
$ws.Range('E35').Value = '  +6.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7243'
$ws.Range('E36').Value = '  +7.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.684'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.869'
$ws.Range('E38').Value = '  +8.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.333'
$ws.Range('E39').Value = '  +8.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9764'
$ws.Range('E40').Value = '  +1.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01708'
$ws.Range('E41').Value = '  +6.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.190'
$ws.Range('E42').Value = '  +4.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4359'
$ws.Range('E43').Value = '  +5.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.42'
$ws.Range('E44').Value = '  +5.32%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.710'
$ws.Range('E46').Value = '  +6.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1283'
$ws.Range('E47').Value = '  +4.72%  '
$ws.Range('E48').Value = '  +4.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.544'
$ws.Range('E49').Value = '  +4.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.37'
$ws.Range('E50').Value = '  +6.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3865'
$ws.Range('E51').Value = '  +7.24%  '
